$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "26.198.70", "21.40") that must
# stay literal text (matching the inlineStr cells in the source workbook).
# Force text number-format before assigning, then clear the format again so
# no extra style index is introduced, only the cell value changes.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.198.70'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.84%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.678.38'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.32%  '
$ws.Range("E4").Value = '  -0.70%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.46'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -3.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5253'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -4.59%  '
$ws.Range("E7").Value = '  -0.67%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2658'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06302'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.40'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07551'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.86%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.690.08'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.63%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.457'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5653'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.29%  '
$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '66.93'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.70%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008041'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -4.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.263.71'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.89%  '
$ws.Range("E18").Value = '  -0.69%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.837'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '188.11'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.44'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.205'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.97%  '
$ws.Range("E23").Value = '  -0.68%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '149.51'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1251'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -5.76%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.595'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -4.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.04'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06181'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.95%  '
$ws.Range("E29").Value = '  -1.55%  '
$ws.Range("E30").Value = '  -3.70%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.496'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -3.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.443'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.637'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.003'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.99%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6067'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.405'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.31%  '
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.101'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.25%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01614'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.083.74'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8688'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.92%  '
$ws.Range("E42").Value = '  -1.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.11'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.828.87'
$ws.Range("D44").ClearFormats()
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.34'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9980'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.88%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.003'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05241'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4257'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.979'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.61%  '
